$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New matchup rows for fall 23 week 12 inputs (A:D), appended after the
# existing data which ends at row 2416.
$newData = @(
    @(3,2,4,1),
    @(6,1,4,2),
    @(3,1,4,2),
    @(6,0,3,2),
    @(3,2,5,0),
    @(5,2,4,0),
    @(4,2,4,1),
    @(5,0,6,3),
    @(2,2,3,1),
    @(7,0,5,2),
    @(5,2,5,0),
    @(2,2,2,1),
    @(5,0,4,3),
    @(2,2,4,0),
    @(3,0,6,3),
    @(4,2,4,0),
    @(3,1,2,2),
    @(5,0,4,2),
    @(6,3,6,0),
    @(4,1,3,2),
    @(5,0,3,2),
    @(4,1,5,2),
    @(4,0,6,2),
    @(4,0,4,2),
    @(4,1,3,2),
    @(4,0,3,2),
    @(5,0,5,3),
    @(3,2,3,1),
    @(4,1,3,2),
    @(5,0,5,2),
    @(6,2,5,0),
    @(6,2,5,0),
    @(3,0,3,3),
    @(3,2,4,1),
    @(3,2,5,1),
    @(3,3,2,0),
    @(5,0,3,2),
    @(4,2,5,0),
    @(7,2,7,0),
    @(3,2,3,1),
    @(4,3,3,0),
    @(5,2,6,0),
    @(4,3,4,0),
    @(3,1,4,2),
    @(5,2,6,0)
)

$startRow = 2417

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $values = $newData[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

$endRow = $startRow + $newData.Count - 1
$lastRow = $endRow
$selCell = "A" + ($lastRow + 1)

$ws.Application.ActiveWindow.ScrollRow = $lastRow - 18
$ws.Range($selCell).Select()
